$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the long "Requirement" descriptions in the same order they were
# originally authored, so the rebuilt shared-strings table lines up with
# the target workbook.
$ws.Range("B3").Value = "사용자는 로그인, 로그아웃 가능"
$ws.Range("B2").Value = "회원은 ID, 비밀번호, 전화번호, 결제 수단, 선호 자전거 유형을 입력하여 회원가입 가능"
$ws.Range("B4").Value = "사용자는 탈퇴 가능"
$ws.Range("B8").Value = "회원은 요금조회 화면에서 대여시간 및 요금을 확인 가능"
$ws.Range("C5").Value = " 대여중인 자전거 리스트 조회(+자전거 반납 +식당 추천)"
$ws.Range("B5").Value = "현재 대여중인 자전거 대여 정보 조회. 대여 정보 조회 리스트에서 특정 자전거를 지정된 대여소에 반납 가능. 반납시 원하면 식당 추천 및 예약 외부 서비스와 연결."
$ws.Range("B5").RowHeight = 42.75

# Rows 6 (자전거 반납) and 7 (식당 추천 여부) are folded into row 5's
# merged description above, so remove them; row 8 (요금 조회, already
# updated) shifts up to become the new row 6.
$ws.Range("A6:C7").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

$ws.Range("A6").Value = 5

# Rows 2 and 6 now hold single-line text, so their old custom height is no
# longer needed; re-fit them back down to the sheet's default row height.
$ws.Rows("2:2").AutoFit()
$ws.Rows("6:6").AutoFit()

$ws.Range("B8").Select()
